$wb = $excel.ActiveWorkbook

# --- Summary sheet: selection D5 -> G22 ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Select()
$wsSummary.Range("G22").Select()

# --- Repayment schedule sheet: selection F6 -> E17, cell O2 -> P2 (shifted one column right) ---
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Range("O2").Copy($wsSchedule.Range("P2"))
$wsSchedule.Range("O2").Clear()
$wsSchedule.Select()
$wsSchedule.Range("E17").Select()

# --- Transactions sheet: values A2 6348 -> 5, A3 691 -> 1, selection D3 -> A2:L3 ---
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("A2").Value = 5
$wsTransactions.Range("A3").Value = 1
$wsTransactions.Select()
$wsTransactions.Range("A2:L3").Select()
